$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.006.07"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "'1.907.38"
$ws.Range("E3").Value = "  -3.17%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D5").Value = "'324.52"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("D8").Value = "'0.3825"
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").Value = "'0.07720"
$ws.Range("E9").Value = "  -2.65%  "
$ws.Range("D10").Value = "'0.9804"
$ws.Range("E10").Value = "  -0.80%  "
$ws.Range("D11").Value = "'22.05"
$ws.Range("E11").Value = "  -3.24%  "
$ws.Range("D12").Value = "'1.918.69"
$ws.Range("E12").Value = "  -5.15%  "
$ws.Range("D13").Value = "'5.667"
$ws.Range("E13").Value = "  -2.33%  "
$ws.Range("D14").Value = "'6.935"
$ws.Range("E14").Value = "  -3.52%  "
$ws.Range("D15").Value = "'0.07033"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("E16").Value = "  -0.08%  "
$ws.Range("D17").Value = "'83.76"
$ws.Range("E17").Value = "  -4.57%  "
$ws.Range("D18").Value = "'0.000009474"
$ws.Range("E18").Value = "  -4.60%  "
$ws.Range("E19").Value = "  -3.35%  "
$ws.Range("D20").Value = "'1.002"
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("D21").Value = "'28.974.52"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").Value = "'5.315"
$ws.Range("E22").Value = "  -4.00%  "
$ws.Range("D23").Value = "'10.87"
$ws.Range("E23").Value = "  -2.48%  "
$ws.Range("D24").Value = "'2.145.66"
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("D25").Value = "'2.095"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "'158.08"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").Value = "'19.04"
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("D28").Value = "'5.653"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").Value = "'117.35"
$ws.Range("D30").Value = "'1.853"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("D31").Value = "'0.09273"
$ws.Range("E31").Value = "  -1.61%  "
$ws.Range("D32").Value = "'0.8613"
$ws.Range("E32").Value = "  -3.37%  "
$ws.Range("E34").Value = "  -5.85%  "
$ws.Range("D35").Value = "'2.953"
$ws.Range("E35").Value = "  -6.89%  "
$ws.Range("D36").Value = "'0.05726"
$ws.Range("E36").Value = "  -1.67%  "
$ws.Range("D37").Value = "'1.152"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("E38").Value = "  -0.02%  "
$ws.Range("D39").Value = "'0.02038"
$ws.Range("E39").Value = "  -3.22%  "
$ws.Range("D40").Value = "'0.5504"
$ws.Range("E40").Value = "  -3.59%  "
$ws.Range("D41").Value = "'7.400"
$ws.Range("E41").Value = "  -4.53%  "
$ws.Range("D42").Value = "'0.1757"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").Value = "'2.827"
$ws.Range("E43").Value = "  +3.36%  "
$ws.Range("D44").Value = "'9.308"
$ws.Range("E44").Value = "  -3.53%  "
$ws.Range("D45").Value = "'0.5181"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").Value = "'11.19"
$ws.Range("E46").Value = "  -5.00%  "
$ws.Range("D47").Value = "'0.06839"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("D48").Value = "'2.054"
$ws.Range("E48").Value = "  -4.86%  "
$ws.Range("D49").Value = "'110.98"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").Value = "'1.775"
$ws.Range("E50").Value = "  -3.05%  "
$ws.Range("D51").Value = "'0.000002550"
$ws.Range("E51").Value = "  -12.30%  "
